$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value that is the same for every
# data row (2..489). It needs to move from serial 45192 to serial 45202.
$lastRow = 489

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}
